$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2022" data column (M) mirroring the existing 2013-2021 table.
# Header cell M4 takes on the same formatting as K4 (s=12).
$ws.Range("M4").Value = 2022
$ws.Range("K4").Copy()
$ws.Range("M4").PasteSpecial(-4122)  # xlPasteFormats

# Data cell M5 (Small enterprises) takes on the same formatting as L5 (s=18).
$ws.Range("M5").Value = 2.2
$ws.Range("L5").Copy()
$ws.Range("M5").PasteSpecial(-4122)  # xlPasteFormats

# Data cell M6 (Medium-sized enterprises) takes on the same formatting as L6 (s=19).
$ws.Range("M6").Value = 1.2
$ws.Range("L6").Copy()
$ws.Range("M6").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# Update the active selection to match the authored state.
$ws.Range("M10").Select()
